# scoretools/test/test1.xlsx — "added readmei and meiutils; read in makemusicodes"
#
# The meifile for the "basecamp" stage row is renamed from the placeholder
# "bc.mei" to the real generated filename, and three new columns record the
# per-ending MEI "name" labels (mc1/mc2/mc3 groups) for the basecamp stage.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# basecamp row (row 2): meifile
$ws.Range("B2").Value = "MSThe Climb (Base Camp).mei"

# basecamp row (row 2): new "name" cells for the three mc groups.
# Set in this order so the shared-string table grows Ending2, Ending1, Ending3.
$ws.Range("U2").Value = "Ending 2"
$ws.Range("O2").Value = "Ending 1"
$ws.Range("Y2").Value = "Ending 3"

# Selection moved from L1 to X4.
$ws.Range("X4").Select() | Out-Null
